# Update "想去人数" (column F) figures to the latest scraped snapshot
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3339
$ws.Range("F4").Value = 31
$ws.Range("F5").Value = 764
$ws.Range("F6").Value = 2386
$ws.Range("F9").Value = 259
$ws.Range("F10").Value = 155
$ws.Range("F11").Value = 377
$ws.Range("F12").Value = 1121
$ws.Range("F13").Value = 473
$ws.Range("F14").Value = 234
$ws.Range("F17").Value = 4955
$ws.Range("F19").Value = 1395
$ws.Range("F20").Value = 3647
$ws.Range("F21").Value = 185
$ws.Range("F22").Value = 223
$ws.Range("F23").Value = 3981
$ws.Range("F24").Value = 5329
$ws.Range("F25").Value = 128
$ws.Range("F27").Value = 583
$ws.Range("F28").Value = 3417
$ws.Range("F29").Value = 402
$ws.Range("F32").Value = 100
$ws.Range("F33").Value = 909
$ws.Range("F34").Value = 1236
$ws.Range("F35").Value = 53
$ws.Range("F37").Value = 1461
$ws.Range("F38").Value = 155
$ws.Range("F39").Value = 1447
$ws.Range("F40").Value = 49
$ws.Range("F42").Value = 931
$ws.Range("F45").Value = 2485
$ws.Range("F46").Value = 93
$ws.Range("F47").Value = 194
$ws.Range("F48").Value = 378

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1043
$ws.Range("F22").Value = 40
$ws.Range("F24").Value = 18

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2786

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2786
$ws.Range("F4").Value = 3339
$ws.Range("F5").Value = 31
$ws.Range("F6").Value = 764
$ws.Range("F7").Value = 2386
$ws.Range("F10").Value = 259
$ws.Range("F11").Value = 1043
$ws.Range("F12").Value = 155
$ws.Range("F13").Value = 377
$ws.Range("F14").Value = 1121
$ws.Range("F15").Value = 474
$ws.Range("F16").Value = 234
$ws.Range("F19").Value = 4955
$ws.Range("F20").Value = 1395
$ws.Range("F21").Value = 3981
$ws.Range("F22").Value = 5329
$ws.Range("F23").Value = 128
$ws.Range("F25").Value = 583
$ws.Range("F26").Value = 3417
$ws.Range("F27").Value = 403
$ws.Range("F30").Value = 100
$ws.Range("F31").Value = 909
$ws.Range("F32").Value = 1236
$ws.Range("F33").Value = 53
$ws.Range("F35").Value = 1461
$ws.Range("F36").Value = 155
$ws.Range("F37").Value = 1448
$ws.Range("F45").Value = 2486
$ws.Range("F46").Value = 93
$ws.Range("F47").Value = 194
$ws.Range("F48").Value = 378

